$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.737.51'
$ws.Range("E2").Value = '  +4.58%  '
$ws.Range("D3").Value = '2.730.03'
$ws.Range("E3").Value = '  +3.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.92'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.27'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +8.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +2.39%  '
$ws.Range("D9").Value = '2.752.15'
$ws.Range("E9").Value = '  +3.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.76'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.112'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.68%  '
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").Value = '3.222.15'
$ws.Range("E14").Value = '  +3.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.48'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.52%  '
$ws.Range("D16").Value = '63.692.45'
$ws.Range("E16").Value = '  +4.56%  '
$ws.Range("E17").Value = '  +7.18%  '
$ws.Range("D18").Value = '2.752.81'
$ws.Range("E18").Value = '  +3.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.04'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.23%  '
$ws.Range("E20").Value = '  +3.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '361.95'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.93'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.996'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.77'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.172'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.63'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("E29").Value = '  +13.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.05'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.18'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.89%  '
$ws.Range("E32").Value = '  +19.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '175.51'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.60'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.07%  '
$ws.Range("E36").Value = '  +7.88%  '
$ws.Range("E37").Value = '  +9.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.84'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +9.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +11.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.27'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '339.14'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.45'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.99'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +13.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.87'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +7.49%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0604'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.37%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.20'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +7.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.647'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0260'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.39'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("E50").Value = '  +2.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.995'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.27%  '
